# Regenerated experiment order (new run timestamps) + sheet reorder.
$wb = $excel.ActiveWorkbook

# Grab references to the worksheets by their current (pre-edit) names,
# before any renaming happens, so the lookups stay unambiguous.
$gng  = $wb.Worksheets.Item("GNG_TO-16512555704071147")
$nb   = $wb.Worksheets.Item("NB_TO-16512555734150155")
$rs   = $wb.Worksheets.Item("RS_TO-16512555734170165")
$tol  = $wb.Worksheets.Item("TOL_TO-16512555734640176")
$vsat = $wb.Worksheets.Item("vSAT_TO-16512555735407321")

# --- Update stimulus-file cell values for the new run --------------------

# GNG task order
$gng.Range("B2").Value = "go_stims-16515890029033554.csv"
$gng.Range("B3").Value = "GNG_stims-165158900291895.csv"
$gng.Range("B4").Value = "go_stims-165158900291895.csv"
$gng.Range("B5").Value = "GNG_stims-16515890029345722.csv"

# NB task order
$nb.Range("B2").Value  = "ZB-match_3-16515890038238697.csv"
$nb.Range("B3").Value  = "TB-16515890073064399.csv"
$nb.Range("B4").Value  = "OB-16515890044374049.csv"
$nb.Range("B5").Value  = "TB-16515890073220308.csv"
$nb.Range("B6").Value  = "TB-1651589006803315.csv"
$nb.Range("B7").Value  = "OB-1651589005347983.csv"
$nb.Range("B8").Value  = "OB-1651589004993494.csv"
$nb.Range("B9").Value  = "ZB-match_8-16515890032560065.csv"
$nb.Range("B10").Value = "ZB-match_7-16515890037113018.csv"

# RS task order (eyes closed / eyes open) is unchanged this run.

# TOL task order
$tol.Range("B2").Value = "MM_stims-16515890029502313.csv"
$tol.Range("B3").Value = "ZM_stims-16515890029345722.csv"
$tol.Range("B4").Value = "MM_stims-16515890029658597.csv"
$tol.Range("B5").Value = "ZM_stims-16515890029502313.csv"
$tol.Range("B6").Value = "MM_stims-16515890029814513.csv"
$tol.Range("B7").Value = "ZM_stims-16515890029658597.csv"

# vSAT task order
$vsat.Range("B2").Value = "SAT_stims-1651589007353281.csv"
$vsat.Range("B3").Value = "SAT_stims-16515890073689065.csv"
$vsat.Range("B4").Value = "vSAT_stims-16515890074001553.csv"
$vsat.Range("B5").Value = "vSAT_stims-1651589007384531.csv"

# --- Rename sheets to reflect the new task-order run ids ------------------
$rs.Name   = "RS_TO-16515890029033554"
$gng.Name  = "GNG_TO-16515890029345722"
$tol.Name  = "TOL_TO-16515890029814513"
$nb.Name   = "NB_TO-1651589007353281"
$vsat.Name = "vSAT_TO-165158900741578"

# --- Reorder tabs: RS, GNG, TOL, NB, vSAT ---------------------------------
# Worksheet handles keep their original index after a rename, so re-fetch a
# fresh reference by name immediately before each Move call, and build the
# final order by repeatedly sending the next sheet to the end of the tab
# strip.
$order = @(
    "RS_TO-16515890029033554",
    "GNG_TO-16515890029345722",
    "TOL_TO-16515890029814513",
    "NB_TO-1651589007353281",
    "vSAT_TO-165158900741578"
)
foreach ($sheetName in $order) {
    $sheet = $wb.Worksheets.Item($sheetName)
    $lastIdx = $wb.Worksheets.Count
    $sheet.Move($null, $wb.Worksheets.Item($lastIdx))
}
